$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to remain Text (they look numeric, e.g. "1.010", "27.003.80")
# so Excel does not silently reinterpret them as numbers and drop formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.003.80'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '1.848.91'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").Value = '1.010'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").Value = '308.98'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").Value = '0.3678'
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("D9").Value = '0.07224'
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = '0.9302'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = '19.76'
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").Value = '0.07738'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '1.903.32'
$ws.Range("E13").Value = '  +3.86%  '
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("D15").Value = '6.439'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '89.04'
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '0.000008649'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D20").Value = '27.022.97'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '14.52'
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("D22").Value = '5.067'
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").Value = '152.66'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").Value = '18.24'
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").Value = '2.015'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").Value = '114.15'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '4.977'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").Value = '0.08852'
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("D31").Value = '3.312'
$ws.Range("E31").Value = '  +4.99%  '
$ws.Range("D32").Value = '1.181'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D33").Value = '0.7411'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '4.500'
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").Value = '2.721'
$ws.Range("E35").Value = '  -4.05%  '
$ws.Range("E36").Value = '  +2.49%  '
$ws.Range("D37").Value = '0.01962'
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("D38").Value = '0.05265'
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '0.5263'
$ws.Range("E40").Value = '  +3.67%  '
$ws.Range("D41").Value = '7.035'
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '8.255'
$ws.Range("E43").Value = '  +1.59%  '
$ws.Range("E44").Value = '  +3.90%  '
$ws.Range("D45").Value = '0.4743'
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("D46").Value = '1.012'
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = '101.84'
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("D49").Value = '65.88'
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("D50").Value = '0.06072'
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").Value = '0.8894'
